$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value parses as a plain number (single decimal
# point) need to be pinned to Text format first, otherwise Excel would silently
# convert them from the inline string they are in the source workbook into a
# numeric cell. Values that keep their "thousands-dot" look (e.g. 26.152.77)
# are not valid numbers so Excel leaves them as text automatically.
$textForceCells = @(
  "D5",
  "D6",
  "D8",
  "D9",
  "D10",
  "D11",
  "D14",
  "D15",
  "D16",
  "D18",
  "D21",
  "D22",
  "D25",
  "D26",
  "D27",
  "D28",
  "D29",
  "D31",
  "D32",
  "D33",
  "D35",
  "D36",
  "D37",
  "D38",
  "D39",
  "D41",
  "D43",
  "D45",
  "D46",
  "D48",
  "D49",
  "D50",
  "D51"
)
foreach ($addr in $textForceCells) {
  $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.152.77"
$ws.Range("E2").Value = "  -2.07%  "
$ws.Range("D3").Value = "1.667.93"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "216.70"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "0.5123"
$ws.Range("E6").Value = "  +4.25%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "0.2637"
$ws.Range("E8").Value = "  +2.23%  "
$ws.Range("D9").Value = "0.06406"
$ws.Range("E9").Value = "  +5.54%  "
$ws.Range("D10").Value = "21.62"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").Value = "0.07419"
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").Value = "1.675.20"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").Value = "0.5814"
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("D15").Value = "0.000008588"
$ws.Range("E15").Value = "  +5.46%  "
$ws.Range("D16").Value = "64.24"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "26.213.43"
$ws.Range("D18").Value = "4.935"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").Value = "189.83"
$ws.Range("E21").Value = "  +4.15%  "
$ws.Range("D22").Value = "6.202"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "7.628"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("D26").Value = "0.1197"
$ws.Range("E26").Value = "  +6.30%  "
$ws.Range("D27").Value = "15.61"
$ws.Range("E27").Value = "  +2.86%  "
$ws.Range("D28").Value = "0.06324"
$ws.Range("E28").Value = "  +14.19%  "
$ws.Range("D29").Value = "1.296"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").Value = "3.531"
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("D32").Value = "3.518"
$ws.Range("E32").Value = "  +2.34%  "
$ws.Range("D33").Value = "1.645"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").Value = "0.6088"
$ws.Range("E35").Value = "  +4.32%  "
$ws.Range("D36").Value = "2.365"
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("D37").Value = "2.653"
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("D38").Value = "6.168"
$ws.Range("E38").Value = "  +5.25%  "
$ws.Range("D39").Value = "0.01605"
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("D40").Value = "1.084.83"
$ws.Range("E40").Value = "  +2.01%  "
$ws.Range("D41").Value = "0.8646"
$ws.Range("E41").Value = "  +1.69%  "
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").Value = "101.04"
$ws.Range("E43").Value = "  +3.29%  "
$ws.Range("D44").Value = "1.816.57"
$ws.Range("E44").Value = "  -1.84%  "
$ws.Range("D45").Value = "0.00000000111"
$ws.Range("E45").Value = "  +5.12%  "
$ws.Range("D46").Value = "56.21"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").Value = "8.101"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").Value = "0.05200"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").Value = "0.4293"
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("D51").Value = "5.895"
$ws.Range("E51").Value = "  +6.14%  "
